# Final draft text & metadata edits before review
# Applies to the "Personnel" worksheet (xl/worksheets/sheet4.xml):
#   - merges/replaces the old "Zoe Sandwith" (row 7) and "Kate Morkeski" (row 8)
#     rows into a single, fuller "Zoe Sandwith" row 7
#   - removes the old row 8 entirely
#   - updates the sheet's selection/view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Remove the old "Kate Morkeski" row (row 8) completely.
$ws.Rows.Item(8).Delete()

# Update row 7 in place with the new personnel record (order matches the
# shared-string insertion order produced by the original authoring edit).
$ws.Range("A7").Value = "Zoe"
$ws.Range("C7").Value = "Sandwith"
$ws.Range("D7").Value = "Northeast U.S. Shelf LTER"
$ws.Range("F7").Value = "0000-0001-9952-9526"
$ws.Range("B7").Value = "O"
$ws.Range("E7").Value = "zoe.sandwith@hakai.org"
$ws.Range("G7").Value = "creator"
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = "OCE-1655687"

# Update the active selection shown when the workbook is reopened.
$ws.Range("J6:J7").Select()

$wb.Save()
